$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.039.20"
$ws.Range("E2").Value = "  +0.69%  "
$ws.Range("D3").Value = "'1.646.06"
$ws.Range("E3").Value = "  +1.16%  "
$ws.Range("E4").Value = "  +0.66%  "
$ws.Range("D5").Value = "'216.84"
$ws.Range("E5").Value = "  +1.10%  "
$ws.Range("E6").Value = "  +1.59%  "
$ws.Range("E7").Value = "  +0.59%  "
$ws.Range("E8").Value = "  +1.00%  "
$ws.Range("E9").Value = "  +1.75%  "
$ws.Range("D10").Value = "'19.70"
$ws.Range("E10").Value = "  +0.59%  "
$ws.Range("D11").Value = "'0.0796"
$ws.Range("E11").Value = "  +1.26%  "
$ws.Range("E12").Value = "  +1.51%  "
$ws.Range("D13").Value = "'1.873.37"
$ws.Range("E13").Value = "  +1.13%  "
$ws.Range("D14").Value = "'1.661.06"
$ws.Range("E14").Value = "  +2.29%  "
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("E16").Value = "  +1.52%  "
$ws.Range("D17").Value = "'63.18"
$ws.Range("E17").Value = "  +0.96%  "
$ws.Range("D18").Value = "'26.017.50"
$ws.Range("E18").Value = "  +0.64%  "
$ws.Range("E19").Value = "  +0.59%  "
$ws.Range("D20").Value = "'193.56"
$ws.Range("E20").Value = "  +0.72%  "
$ws.Range("D21").Value = "'4.37"
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("E23").Value = "  +0.53%  "
$ws.Range("E24").Value = "  +8.20%  "
$ws.Range("E25").Value = "  +2.27%  "
$ws.Range("D26").Value = "'144.62"
$ws.Range("E26").Value = "  +1.48%  "
$ws.Range("E27").Value = "  +0.71%  "
$ws.Range("D28").Value = "'6.94"
$ws.Range("E28").Value = "  +1.28%  "
$ws.Range("D29").Value = "'15.57"
$ws.Range("E29").Value = "  +1.03%  "
$ws.Range("E30").Value = "  +1.19%  "
$ws.Range("E31").Value = "  +0.94%  "
$ws.Range("E32").Value = "  -0.42%  "
$ws.Range("E33").Value = "  +1.63%  "
$ws.Range("D34").Value = "'1.54"
$ws.Range("E34").Value = "  -2.37%  "
$ws.Range("E35").Value = "  +3.07%  "
$ws.Range("E36").Value = "  +0.70%  "
$ws.Range("D37").Value = "'1.134.45"
$ws.Range("E37").Value = "  +0.42%  "
$ws.Range("E38").Value = "  -0.51%  "
$ws.Range("E39").Value = "  +0.30%  "
$ws.Range("D40").Value = "'0.0158"
$ws.Range("E40").Value = "  +1.17%  "
$ws.Range("E41").Value = "  +0.85%  "
$ws.Range("D42").Value = "'99.63"
$ws.Range("E42").Value = "  +0.22%  "
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("D44").Value = "'1.781.98"
$ws.Range("E44").Value = "  +1.06%  "
$ws.Range("E45").Value = "  +3.46%  "
$ws.Range("D46").Value = "'56.86"
$ws.Range("E46").Value = "  +1.45%  "
$ws.Range("E47").Value = "  -0.18%  "
$ws.Range("E48").Value = "  +0.66%  "
$ws.Range("D49").Value = "'7.70"
$ws.Range("E49").Value = "  +1.73%  "
$ws.Range("E50").Value = "  +0.33%  "
$ws.Range("E51").Value = "  +0.36%  "
